$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @"
2,D,Neutrophils
2,G,0.60844
2,H,1.82532
2,I,0.01418783042133501
2,J,0.01418783042133501
2,K,1
2,L,0.3333333333333333
2,M,0.1258873333333333
2,N,0.377662
2,O,0.9923353626973638
2,P,0.9923353626973638
2,Q,0.07659488909333333
2,R,0.6893540018400001
2,S,0.01407908584704416
2,T,0.01407908584704416
3,G,0.60844
3,H,1.82532
3,I,0.01418783042133501
3,J,0.01418783042133501
3,M,0.0009723333333333333
3,N,0.002917
3,O,0.007664637302636247
3,P,0.007664637302636246
3,Q,0.0005916064933333333
3,R,0.00532445844
3,S,0.0001087445742908416
3,T,0.0001087445742908416
4,D,Neutrophils
4,I,0.02409258886165303
4,J,0.02409258886165303
4,K,1
4,L,0.3333333333333333
4,M,0.1258873333333333
4,N,0.377662
4,O,0.9923353626973638
4,P,0.9923353626973638
4,Q,0.1300670445746666
4,R,1.170603401172
4,S,0.02390792790634693
4,T,0.02390792790634692
5,I,0.02409258886165303
5,J,0.02409258886165303
5,M,0.0009723333333333333
5,N,0.002917
5,O,0.007664637302636247
5,P,0.007664637302636246
5,Q,0.001004616744666666
5,R,0.009041550701999999
5,S,0.0001846609553061044
5,T,0.0001846609553061043
6,D,Neutrophils
6,G,6.042074333333333
6,H,18.126223
6,I,0.1408913385616233
6,J,0.1408913385616233
6,K,1
6,L,0.3333333333333333
6,M,0.1258873333333333
6,N,0.377662
6,O,0.9923353626973638
6,P,0.9923353626973638
6,Q,0.760620625625111
6,R,6.845585630626
6,S,0.1398114575524655
6,T,0.1398114575524655
7,G,6.042074333333333
7,H,18.126223
7,I,0.1408913385616233
7,J,0.1408913385616233
7,M,0.0009723333333333333
7,N,0.002917
7,O,0.007664637302636247
7,P,0.007664637302636246
7,Q,0.005874910276777777
7,R,0.052874192491
7,S,0.001079881009157771
7,T,0.001079881009157771
8,D,Neutrophils
8,G,0.4072233333333333
8,H,1.22167
8,I,0.009495785281940885
8,J,0.009495785281940885
8,K,1
8,L,0.3333333333333333
8,M,0.1258873333333333
8,N,0.377662
8,O,0.9923353626973638
8,P,0.9923353626973638
8,Q,0.05126425950444444
8,R,0.46137833554
8,S,0.009423003531851097
8,T,0.009423003531851097
9,G,0.4072233333333333
9,H,1.22167
9,I,0.009495785281940885
9,J,0.009495785281940885
9,M,0.0009723333333333333
9,N,0.002917
9,O,0.007664637302636247
9,P,0.007664637302636246
9,Q,0.0003959568211111111
9,R,0.00356361139
9,S,0.00007278175008978836
9,T,0.00007278175008978836
10,D,Neutrophils
10,G,31.00247266666667
10,H,93.007418
10,I,0.72292719879814
10,J,0.72292719879814
10,K,1
10,L,0.3333333333333333
10,M,0.1258873333333333
10,N,0.377662
10,O,0.9923353626973638
10,P,0.9923353626973638
10,Q,3.902818610746222
10,R,35.125367496716
10,S,0.7173862240231416
10,T,0.7173862240231416
11,G,31.00247266666667
11,H,93.007418
11,I,0.72292719879814
11,J,0.72292719879814
11,M,0.0009723333333333333
11,N,0.002917
11,O,0.007664637302636247
11,P,0.007664637302636246
11,Q,0.03014473758955555
11,R,0.271302638306
11,S,0.005540974774998554
11,T,0.005540974774998554
12,D,Neutrophils
12,G,3.791227666666666
12,H,11.373683
12,I,0.08840525807530777
12,J,0.08840525807530777
12,K,1
12,L,0.3333333333333333
12,M,0.1258873333333333
12,N,0.377662
12,O,0.9923353626973638
12,P,0.9923353626973638
12,Q,0.4772675410162222
12,R,4.295407869146
12,S,0.08772766383651459
12,T,0.08772766383651459
13,G,3.791227666666666
13,H,11.373683
13,I,0.08840525807530777
13,J,0.08840525807530777
13,M,0.0009723333333333333
13,N,0.002917
13,O,0.007664637302636247
13,P,0.007664637302636246
13,Q,0.003686337034555555
13,R,0.033177033311
13,S,0.0006775942387931882
13,T,0.0006775942387931882
"@

$rows = $changes -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $r = $parts[0]
    $c = $parts[1]
    $val = $parts[2]
    $ws.Range("$c$r").Value = $val
}

Write-Host "Applied $($rows.Count) cell updates"
